# Auto-generated: apply crypto price/volume updates for Wed Oct  4 04:41:37 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.420.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("E6").Value = "  +3.51%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.87%  "
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  -2.22%  "
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.872.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("E16").Value = "  -2.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.379.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -5.24%  "
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("E32").Value = "  -2.14%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.412.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.561"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.59%  "
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.788"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.782.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.43%  "
$ws.Range("E48").Value = "  -4.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0104"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("E51").Value = "  -3.64%  "
